$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Move/insert rows to make room for the two new "my static" rows,
#     and for extra blank spacing above the helper data block lower down.
$ws.Rows("8:9").Insert()
$ws.Rows("15:20").Insert()

# --- New data rows: "my static (1)" / chunk size = 1, and "my static (MAX)" / max chunk size.
#     Shared-string order matters: write B8, B9 first (my static (1), my static (MAX)),
#     then A8/A10 (chunk size = 1), then A9 (max chunk size) so the new shared strings
#     are appended to sharedStrings.xml in that exact order.
$ws.Range("B8").Value = "my static (1)"
$ws.Range("D8").Value = 6.3017
$ws.Range("E8").Value = 6.0128
$ws.Range("F8").Value = 6.5102

$ws.Range("B9").Value = "my static (MAX)"
$ws.Range("D9").Value = 5.9861
$ws.Range("E9").Value = 5.6196
$ws.Range("F9").Value = 5.6618

$ws.Range("A8").Value = "chunk size = 1"
$ws.Range("A9").Value = "max chunk size"
$ws.Range("A10").Value = "chunk size = 1"

# --- Updated measurement for the "static" row (now row 10) in the last column.
$ws.Range("F10").Value = 8.235

# --- Widen the new label columns (A:B) slightly (target stored width 14.75).
$ws.Columns("A:B").ColumnWidth = 14.035714285714286

# --- Grow table4 (this sheet's table) down to the new "guided" row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B6:F12"))

# --- Selection / window bookkeeping to match the saved view.
$ws.Range("C1").Select() | Out-Null
$win = $wb.Windows.Item(1)
$win.Left = 2790
